# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "42.993.38"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.294.23"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "300.28"
$ws.Range("E5").Value = "  +0.22%  "
Set-TextValue "D6" "99.13"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.50%  "
Set-TextValue "D10" "36.21"
$ws.Range("E10").Value = "  +8.47%  "
$ws.Range("E12").Value = "  +1.06%  "
Set-TextValue "D13" "18.16"
$ws.Range("E13").Value = "  +8.32%  "
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "2.651.90"
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.274.64"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "0.798"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "42.881.22"
$ws.Range("E18").Value = "  +0.26%  "
Set-TextValue "D19" "12.52"
$ws.Range("E19").Value = "  +8.62%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +1.24%  "
Set-TextValue "D22" "67.71"
$ws.Range("E22").Value = "  +0.40%  "
Set-TextValue "D23" "235.58"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +10.32%  "
$ws.Range("E25").Value = "  +0.41%  "
Set-TextValue "D26" "2.44"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  +2.73%  "
Set-TextValue "D29" "34.46"
$ws.Range("E29").Value = "  +2.61%  "
Set-TextValue "D30" "166.93"
$ws.Range("E30").Value = "  -0.07%  "
Set-TextValue "D31" "9.11"
$ws.Range("E31").Value = "  +0.01%  "
Set-TextValue "D32" "0.999"
$ws.Range("E32").Value = "  -0.06%  "
Set-TextValue "D33" "5.01"
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  +1.19%  "
Set-TextValue "D37" "0.0686"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  -0.13%  "
Set-TextValue "D42" "2.30"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").Value = "1.964.22"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("E46").Value = "  +1.86%  "
Set-TextValue "D47" "17.42"
$ws.Range("E47").Value = "  -0.74%  "
Set-TextValue "D48" "55.42"
$ws.Range("E48").Value = "  +4.78%  "
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").Value = "2.518.60"
Set-TextValue "D51" "70.66"
$ws.Range("E51").Value = "  +1.14%  "
